$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, applied as text to preserve the
# original text-formatted (non-numeric) storage of these cells.
$updates = [ordered]@{
    "D2" = '27.950.49'
    "E2" = '  -2.01%  '
    "D3" = '1.868.44'
    "E3" = '  -2.61%  '
    "D4" = '1.001'
    "E4" = '  -0.11%  '
    "D5" = '312.32'
    "E5" = '  -1.23%  '
    "D6" = '1.001'
    "E6" = '  -0.11%  '
    "D7" = '0.4993'
    "E7" = '  -2.35%  '
    "D8" = '0.3827'
    "E8" = '  -3.55%  '
    "D9" = '0.08914'
    "E9" = '  -8.32%  '
    "D10" = '1.119'
    "D11" = '41.47'
    "E11" = '  -1.47%  '
    "D12" = '6.364'
    "E12" = '  -1.66%  '
    "D13" = '20.71'
    "E13" = '  -1.39%  '
    "D14" = '1.859.71'
    "E14" = '  -3.19%  '
    "D15" = '7.232'
    "E15" = '  -2.35%  '
    "D16" = '1.002'
    "E16" = '  -0.07%  '
    "D17" = '0.00001098'
    "E17" = '  -2.92%  '
    "D18" = '91.07'
    "E18" = '  -3.13%  '
    "D19" = '0.06672'
    "E19" = '  +0.03%  '
    "D20" = '18.01'
    "E20" = '  -0.64%  '
    "E21" = '  -0.05%  '
    "D22" = '6.113'
    "E22" = '  -2.67%  '
    "D23" = '28.006.24'
    "E23" = '  -2.05%  '
    "D24" = '11.49'
    "E24" = '  +0.12%  '
    "D25" = '2.277'
    "E25" = '  -1.51%  '
    "B26" = 'LEO'
    "C26" = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    "D26" = '3.399'
    "E26" = '  +0.68%  '
    "B27" = 'WrappedliquidstakedEther2.0'
    "C27" = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    "D27" = '2.071.70'
    "E27" = '  -3.15%  '
    "B28" = 'LidoDAOToken'
    "C28" = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    "D28" = '2.507'
    "E28" = '  -6.26%  '
    "B29" = 'Monero'
    "C29" = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    "D29" = '158.20'
    "E29" = '  +0.01%  '
    "B30" = 'EthereumClassic'
    "C30" = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    "D30" = '20.67'
    "E30" = '  -2.52%  '
    "B31" = 'BitcoinCash'
    "C31" = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    "D31" = '126.11'
    "E31" = '  -2.02%  '
    "B32" = 'Stellar'
    "C32" = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    "D32" = '0.1058'
    "E32" = '  -0.96%  '
    "B33" = 'ImmutableX'
    "C33" = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    "D33" = '1.054'
    "E33" = '  -4.71%  '
    "B34" = 'Filecoin'
    "C34" = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    "D34" = '5.587'
    "E34" = '  -1.95%  '
    "B35" = 'HuobiToken'
    "C35" = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    "D35" = '3.586'
    "E35" = '  -1.40%  '
    "B36" = 'FraxShare'
    "C36" = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    "D36" = '9.392'
    "E36" = '  -4.06%  '
    "B37" = 'Hedera'
    "C37" = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    "D37" = '0.06551'
    "E37" = '  -2.33%  '
    "B38" = 'VeChain'
    "C38" = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    "D38" = '0.02407'
    "E38" = '  -1.41%  '
    "B39" = 'Algorand'
    "C39" = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    "D39" = '0.2187'
    "E39" = '  -1.50%  '
    "B40" = 'TrustWalletToken'
    "C40" = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    "D40" = '1.279'
    "E40" = '  +5.88%  '
    "B41" = 'ARBITRUM'
    "C41" = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    "D41" = '1.199'
    "E41" = '  -4.70%  '
    "B42" = 'Aptos'
    "C42" = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    "D42" = '11.55'
    "E42" = '  -0.89%  '
    "B43" = 'TheSandbox'
    "C43" = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    "D43" = '0.6375'
    "E43" = '  -1.05%  '
    "B44" = 'InternetComputer(DFINITY)'
    "C44" = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    "D44" = '4.906'
    "E44" = '  -2.95%  '
    "B45" = 'Frax'
    "C45" = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    "D45" = '1.000'
    "E45" = '  -0.10%  '
    "D46" = '13.26'
    "E46" = '  -3.26%  '
    "B47" = 'Decentraland'
    "C47" = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    "D47" = '0.6002'
    "E47" = '  -1.19%  '
    "B48" = 'WEMIXTOKEN'
    "C48" = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    "D48" = '1.281'
    "E48" = '  -0.26%  '
    "B49" = 'PancakeSwap'
    "C49" = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    "D49" = '3.673'
    "E49" = '  -2.82%  '
    "B50" = 'NEARProtocol'
    "C50" = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    "D50" = '1.993'
    "E50" = '  -3.30%  '
    "B51" = 'EOS'
    "C51" = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
    "D51" = '1.223'
    "E51" = '  +2.19%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text number format so Excel does not reinterpret strings such
    # as "1.001" or "0.00001098" as numbers, matching the source file where
    # these columns are stored as plain text.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    # Drop back to the default "Normal" style so no stray style index is
    # left behind on cells that originally had no explicit style.
    $cell.Style = "Normal"
}
